# Insert a new weekly price record as row 31, pushing the existing
# rows 31-73 down to 32-74 (dimension grows from A1:R73 to A1:R74).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("31:31").Insert()

$ws.Cells.Item(31, 1).Value = 1
$ws.Cells.Item(31, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(31, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(31, 4).Value = 44935
$ws.Cells.Item(31, 5).Value = 15
$ws.Cells.Item(31, 6).Value = 100112027
$ws.Cells.Item(31, 7).Value = "Melón"
$ws.Cells.Item(31, 8).Value = "Calameño"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 130
$ws.Cells.Item(31, 11).Value = 9000
$ws.Cells.Item(31, 12).Value = 10000
$ws.Cells.Item(31, 13).Value = 9615
$ws.Cells.Item(31, 14).Value = '$/caja 18 unidades'
$ws.Cells.Item(31, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(31, 16).Value = 534
$ws.Cells.Item(31, 17).Value = 18
$ws.Cells.Item(31, 18).Value = "Hortaliza"
